# Generate Report for handback
# Adds a new handback-status row for file "a4807e53-9cc8-4008-b095-bb857e88121d"
# to the Overview / zh-cn / de-de worksheets, mirroring the existing rows for
# 043ae17f-bee3-45e7-8ada-d4ceb040a5d1 and 2bcadd7e-a878-42d0-8d49-4877ee27b325.

$wb = $excel.ActiveWorkbook

$newId   = "a4807e53-9cc8-4008-b095-bb857e88121d"
$mdName  = "$newId.md"
$zhXlf   = "$newId.fc4523630b7ff9b7c970a156bc05c0e68234f618.zh-cn.xlf"
$deXlf   = "$newId.fc4523630b7ff9b7c970a156bc05c0e68234f618.de-de.xlf"

$inSync  = "Handed back: in sync with en-US"
$include = "Include"

$zhHandoffDt  = "2016-01-18 12:15:57"
$zhHandbackDt = "2016-01-18 12:16:37"
$deHandoffDt  = "2016-01-18 12:16:08"
$deHandbackDt = "2016-01-18 12:16:53"

# Helper colour/underline to mirror the workbook's existing "hyperlink" look
# (underline + cornflower-blue font) so the new cells read consistently with
# the rest of the sheet.
$hlUnderline = 2        # xlUnderlineStyleSingle
$hlColor     = 15570276 # OLE (BGR) for RGB 6495ED

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $hlUnderline
    $rng.Font.Color = $hlColor
}

# ---------------------------------------------------------------------------
# Sheet "Overview" -> row 4
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$a4 = $wsOverview.Range("A4")
$a4.Value = $mdName
$wsOverview.Hyperlinks.Add($a4, "https://github.com/OpenLocalizationTest/oltest/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $a4

$wsOverview.Range("B4").Value = $inSync
$wsOverview.Range("C4").Value = $inSync

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$zhA4 = $wsZh.Range("A4")
$zhA4.Value = $mdName
$wsZh.Hyperlinks.Add($zhA4, "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $zhA4

$wsZh.Range("B4").Value = $inSync

$zhC4 = $wsZh.Range("C4")
$zhC4.Value = $zhXlf
$wsZh.Hyperlinks.Add($zhC4, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf", "", "", $zhXlf)
Style-AsHyperlink $zhC4

$wsZh.Range("D4").Value = $zhHandoffDt

$zhE4 = $wsZh.Range("E4")
$zhE4.Value = $mdName
$wsZh.Hyperlinks.Add($zhE4, "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $zhE4

$zhF4 = $wsZh.Range("F4")
$zhF4.Value = $zhXlf
$wsZh.Hyperlinks.Add($zhF4, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$zhXlf", "", "", $zhXlf)
Style-AsHyperlink $zhF4

$wsZh.Range("G4").Value = $zhHandbackDt
$wsZh.Range("H4").Value = $include

# ---------------------------------------------------------------------------
# Sheet "de-de" -> row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$deA4 = $wsDe.Range("A4")
$deA4.Value = $mdName
$wsDe.Hyperlinks.Add($deA4, "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $deA4

$wsDe.Range("B4").Value = $inSync

$deC4 = $wsDe.Range("C4")
$deC4.Value = $deXlf
$wsDe.Hyperlinks.Add($deC4, "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf", "", "", $deXlf)
Style-AsHyperlink $deC4

$wsDe.Range("D4").Value = $deHandoffDt

$deE4 = $wsDe.Range("E4")
$deE4.Value = $mdName
$wsDe.Hyperlinks.Add($deE4, "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/e2e/$mdName", "", "", $mdName)
Style-AsHyperlink $deE4

$deF4 = $wsDe.Range("F4")
$deF4.Value = $deXlf
$wsDe.Hyperlinks.Add($deF4, "https://github.com/OpenLocalizationTestOrg/olhandback/blob/fc4523630b7ff9b7c970a156bc05c0e68234f618/ol-handback/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$deXlf", "", "", $deXlf)
Style-AsHyperlink $deF4

$wsDe.Range("G4").Value = $deHandbackDt
$wsDe.Range("H4").Value = $include

Write-Host "Added handback row for $newId to Overview, zh-cn, de-de sheets."
